$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-17T07:09:33.705290"
    3 = "2025-10-17T07:09:33.705290"
    4 = "2025-10-17T07:09:33.705290"
    5 = "2025-10-17T07:09:33.705290"
    6 = "2025-10-17T07:09:33.705290"
    7 = "2025-10-17T07:09:33.705290"
    8 = "2025-10-17T07:09:33.705290"
    9 = "2025-10-17T07:09:33.706289"
    10 = "2025-10-17T07:09:33.706289"
    11 = "2025-10-17T07:09:33.706289"
    12 = "2025-10-17T07:09:33.706289"
    13 = "2025-10-17T07:09:33.706289"
    14 = "2025-10-17T07:09:33.706289"
    15 = "2025-10-17T07:09:33.706289"
    16 = "2025-10-17T07:09:33.706289"
    17 = "2025-10-17T07:09:33.706289"
    18 = "2025-10-17T07:09:33.706289"
    19 = "2025-10-17T07:09:33.706289"
    20 = "2025-10-17T07:09:33.707289"
    21 = "2025-10-17T07:09:33.707289"
    22 = "2025-10-17T07:09:33.707289"
    23 = "2025-10-17T07:09:33.707289"
    24 = "2025-10-17T07:09:33.707289"
    25 = "2025-10-17T07:09:33.707289"
    26 = "2025-10-17T07:09:33.707289"
    27 = "2025-10-17T07:09:33.707289"
    28 = "2025-10-17T07:09:33.707289"
    29 = "2025-10-17T07:09:33.707289"
    30 = "2025-10-17T07:09:33.708292"
    31 = "2025-10-17T07:09:33.708292"
    32 = "2025-10-17T07:09:33.709295"
    33 = "2025-10-17T07:09:33.710296"
    34 = "2025-10-17T07:09:33.710296"
    35 = "2025-10-17T07:09:33.710989"
    36 = "2025-10-17T07:09:33.710989"
    37 = "2025-10-17T07:09:33.710989"
    38 = "2025-10-17T07:09:33.710989"
    39 = "2025-10-17T07:09:33.711503"
    40 = "2025-10-17T07:09:33.711503"
    41 = "2025-10-17T07:09:33.711503"
    42 = "2025-10-17T07:09:33.712251"
    43 = "2025-10-17T07:09:33.712251"
    44 = "2025-10-17T07:09:33.712251"
    45 = "2025-10-17T07:09:33.712770"
    46 = "2025-10-17T07:09:33.792320"
    47 = "2025-10-17T07:09:33.792320"
    48 = "2025-10-17T07:09:33.792320"
    49 = "2025-10-17T07:09:33.792320"
    50 = "2025-10-17T07:09:33.792320"
    51 = "2025-10-17T07:09:33.792320"
    52 = "2025-10-17T07:09:33.792320"
    53 = "2025-10-17T07:09:33.792320"
    54 = "2025-10-17T07:09:33.792320"
    55 = "2025-10-17T07:09:33.792320"
    56 = "2025-10-17T07:09:33.792320"
    57 = "2025-10-17T07:09:33.792320"
    58 = "2025-10-17T07:09:33.792320"
    59 = "2025-10-17T07:09:33.792320"
    60 = "2025-10-17T07:09:33.792320"
    61 = "2025-10-17T07:09:33.792320"
    62 = "2025-10-17T07:09:33.792320"
    63 = "2025-10-17T07:09:33.792320"
    64 = "2025-10-17T07:09:33.792320"
    65 = "2025-10-17T07:09:33.792320"
    66 = "2025-10-17T07:09:33.792320"
    67 = "2025-10-17T07:09:33.792320"
    68 = "2025-10-17T07:09:33.792320"
    69 = "2025-10-17T07:09:33.792320"
    70 = "2025-10-17T07:09:33.792320"
    71 = "2025-10-17T07:09:33.792320"
    72 = "2025-10-17T07:09:33.792320"
    73 = "2025-10-17T07:09:33.792320"
    74 = "2025-10-17T07:09:33.792320"
    75 = "2025-10-17T07:09:33.880392"
    76 = "2025-10-17T07:09:33.880392"
    77 = "2025-10-17T07:09:33.880392"
    78 = "2025-10-17T07:09:33.880392"
    79 = "2025-10-17T07:09:33.880392"
    80 = "2025-10-17T07:09:33.880392"
    81 = "2025-10-17T07:09:33.880392"
    82 = "2025-10-17T07:09:33.880392"
    83 = "2025-10-17T07:09:33.880392"
    84 = "2025-10-17T07:09:33.880392"
    85 = "2025-10-17T07:09:33.880392"
    86 = "2025-10-17T07:09:33.880392"
    87 = "2025-10-17T07:09:33.880392"
    88 = "2025-10-17T07:09:33.880392"
    89 = "2025-10-17T07:09:33.880392"
    90 = "2025-10-17T07:09:33.880392"
    91 = "2025-10-17T07:09:33.880392"
    92 = "2025-10-17T07:09:33.880392"
    93 = "2025-10-17T07:09:33.880392"
    94 = "2025-10-17T07:09:33.880392"
    95 = "2025-10-17T07:09:33.880392"
    96 = "2025-10-17T07:09:33.880392"
    97 = "2025-10-17T07:09:33.880392"
    98 = "2025-10-17T07:09:33.880392"
    99 = "2025-10-17T07:09:33.880392"
    100 = "2025-10-17T07:09:33.880392"
    101 = "2025-10-17T07:09:33.880392"
    102 = "2025-10-17T07:09:33.880392"
    103 = "2025-10-17T07:09:33.991598"
    104 = "2025-10-17T07:09:33.993631"
    105 = "2025-10-17T07:09:33.993631"
    106 = "2025-10-17T07:09:33.994146"
    107 = "2025-10-17T07:09:33.994146"
    108 = "2025-10-17T07:09:33.994146"
    109 = "2025-10-17T07:09:33.994146"
    110 = "2025-10-17T07:09:33.994146"
    111 = "2025-10-17T07:09:33.994146"
    112 = "2025-10-17T07:09:33.994146"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}

Write-Host "Updated $($timestamps.Count) timestamp cells"
